# Auto-generated Excel COM-interop script to apply profit recalculation updates
# across multiple worksheets of the Tiamat_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 272.8
$ws.Range("J18").Value = 1200
$ws.Range("L18").Value = 1200
$ws.Range("N18").Value = -1768

$ws.Range("H39").Value = 1314.55
$ws.Range("I39").Value = 85.888885
$ws.Range("J39").Value = 2319.818
$ws.Range("K39").Value = 257.666655
$ws.Range("L39").Value = 6959.454000000001
$ws.Range("M39").Value = 38.33334500000001
$ws.Range("N39").Value = -7551.454000000001

$ws.Range("H58").Value = 2812.6667
$ws.Range("I58").Value = 231.66667
$ws.Range("J58").Value = 4533.3335
$ws.Range("K58").Value = 695.00001
$ws.Range("L58").Value = 13600.0005
$ws.Range("M58").Value = -545.00001
$ws.Range("N58").Value = -13900.0005

$ws.Range("H64").Value = 73998.14
$ws.Range("I64").Value = 2997.8333
$ws.Range("K64").Value = 2997.8333
$ws.Range("M64").Value = -2749.8333

$ws.Range("H67").Value = 73998.14
$ws.Range("I67").Value = 2997.8333
$ws.Range("K67").Value = 2997.8333
$ws.Range("M67").Value = -2139.8333

$ws.Range("H74").Value = 3542
$ws.Range("I74").Value = 3551.0715
$ws.Range("J74").Value = 3499.6667
$ws.Range("K74").Value = 3551.0715
$ws.Range("L74").Value = 3499.6667
$ws.Range("M74").Value = -2615.0715
$ws.Range("N74").Value = -5371.6667

$ws.Range("H76").Value = 37040130
$ws.Range("I76").Value = 41669784
$ws.Range("J76").Value = 2888.6667
$ws.Range("K76").Value = 41669784
$ws.Range("L76").Value = 2888.6667
$ws.Range("M76").Value = -41669469
$ws.Range("N76").Value = -3518.6667

$ws.Range("H77").Value = 3542
$ws.Range("I77").Value = 3551.0715
$ws.Range("J77").Value = 3499.6667
$ws.Range("K77").Value = 17755.3575
$ws.Range("L77").Value = 17498.3335
$ws.Range("M77").Value = -13075.3575
$ws.Range("N77").Value = -26858.3335

$ws.Range("H79").Value = 37040130
$ws.Range("I79").Value = 41669784
$ws.Range("J79").Value = 2888.6667
$ws.Range("K79").Value = 41669784
$ws.Range("L79").Value = 2888.6667
$ws.Range("M79").Value = -41668692
$ws.Range("N79").Value = -5072.6667

$ws.Range("H100").Value = 19401.2
$ws.Range("I100").Value = 19000
$ws.Range("J100").Value = 19501.5
$ws.Range("K100").Value = 19000
$ws.Range("L100").Value = 19501.5
$ws.Range("M100").Value = -18459
$ws.Range("N100").Value = -20583.5

$ws.Range("H132").Value = 162079.27
$ws.Range("I132").Value = 3367.6785
$ws.Range("J132").Value = 1431772
$ws.Range("K132").Value = 10103.0355
$ws.Range("L132").Value = 4295316
$ws.Range("M132").Value = -7573.0355
$ws.Range("N132").Value = -4300376


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16732.57
$ws.Range("I32").Value = 10712.2
$ws.Range("J32").Value = 34793.68
$ws.Range("K32").Value = 10712.2
$ws.Range("L32").Value = 34793.68
$ws.Range("M32").Value = -10425.2
$ws.Range("N32").Value = -35367.68

$ws.Range("H92").Value = 31506.25
$ws.Range("J92").Value = 31506.25
$ws.Range("L92").Value = 31506.25
$ws.Range("N92").Value = -36498.25

$ws.Range("H110").Value = 1271.2632
$ws.Range("I110").Value = 645.125
$ws.Range("J110").Value = 1726.6364
$ws.Range("K110").Value = 645.125
$ws.Range("L110").Value = 1726.6364
$ws.Range("M110").Value = 1399.875
$ws.Range("N110").Value = -5816.6364

$ws.Range("H125").Value = 35714.285
$ws.Range("J125").Value = 35714.285
$ws.Range("L125").Value = 35714.285
$ws.Range("N125").Value = -45554.285


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1555.7646
$ws.Range("I20").Value = 1357.375
$ws.Range("J20").Value = 1732.1111
$ws.Range("K20").Value = 1357.375
$ws.Range("L20").Value = 1732.1111
$ws.Range("M20").Value = -1110.375
$ws.Range("N20").Value = -2226.1111

$ws.Range("H86").Value = 335100.9
$ws.Range("I86").Value = 1442.5
$ws.Range("J86").Value = 779978.75
$ws.Range("K86").Value = 1442.5
$ws.Range("L86").Value = 779978.75
$ws.Range("M86").Value = -319.5
$ws.Range("N86").Value = -782224.75

$ws.Range("H89").Value = 335100.9
$ws.Range("I89").Value = 1442.5
$ws.Range("J89").Value = 779978.75
$ws.Range("K89").Value = 7212.5
$ws.Range("L89").Value = 3899893.75
$ws.Range("M89").Value = -1596.5
$ws.Range("N89").Value = -3911125.75

$ws.Range("H134").Value = 41710896
$ws.Range("I134").Value = 2565.75
$ws.Range("J134").Value = 125127550
$ws.Range("K134").Value = 7697.25
$ws.Range("L134").Value = 375382650
$ws.Range("M134").Value = -5162.25
$ws.Range("N134").Value = -375387720


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34315.742
$ws.Range("I31").Value = 42110.88
$ws.Range("J31").Value = 20395.857
$ws.Range("K31").Value = 42110.88
$ws.Range("L31").Value = 20395.857
$ws.Range("M31").Value = -41815.88
$ws.Range("N31").Value = -20985.857

$ws.Range("H34").Value = 34315.742
$ws.Range("I34").Value = 42110.88
$ws.Range("J34").Value = 20395.857
$ws.Range("K34").Value = 42110.88
$ws.Range("L34").Value = 20395.857
$ws.Range("M34").Value = -41908.88
$ws.Range("N34").Value = -20799.857


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 795.95
$ws.Range("I131").Value = 452.5
$ws.Range("J131").Value = 825.81525
$ws.Range("K131").Value = 1357.5
$ws.Range("L131").Value = 2477.44575
$ws.Range("M131").Value = 3682.5
$ws.Range("N131").Value = -12557.44575


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1665.186
$ws.Range("I61").Value = 1455.9656
$ws.Range("J61").Value = 2098.5715
$ws.Range("K61").Value = 1455.9656
$ws.Range("L61").Value = 2098.5715
$ws.Range("M61").Value = -1253.9656
$ws.Range("N61").Value = -2502.5715

$ws.Range("H113").Value = 1665.186
$ws.Range("I113").Value = 1455.9656
$ws.Range("J113").Value = 2098.5715
$ws.Range("K113").Value = 1455.9656
$ws.Range("L113").Value = 2098.5715
$ws.Range("M113").Value = 714.0344
$ws.Range("N113").Value = -6438.5715


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 19977
$ws.Range("J54").Value = 19977
$ws.Range("L54").Value = 19977
$ws.Range("N54").Value = -21017

$ws.Range("H62").Value = 5050.1113
$ws.Range("I62").Value = 4722.4443
$ws.Range("J62").Value = 5377.778
$ws.Range("K62").Value = 4722.4443
$ws.Range("L62").Value = 5377.778
$ws.Range("M62").Value = -4098.4443
$ws.Range("N62").Value = -6625.778

$ws.Range("H65").Value = 5050.1113
$ws.Range("I65").Value = 4722.4443
$ws.Range("J65").Value = 5377.778
$ws.Range("K65").Value = 23612.2215
$ws.Range("L65").Value = 26888.89
$ws.Range("M65").Value = -20492.2215
$ws.Range("N65").Value = -33128.89

